$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3: field name "submission_date" -> "date_submitted"
$ws.Range("B3").Value = "date_submitted"

# Update row 3 description: append new sentence about date format
$ws.Range("D3").Value = "Date report was received by CTP; this is the earliest date of report receipt, either to Safety Reporting Portal (SRP) or by other means. Date follows format: ``YYYYmmdd``."

# Move the selection/active cell to D3 (reflects the cell the author last edited)
$ws.Range("D3").Select()

# Row 2 now wraps to an explicit height (matches target revision)
$ws.Rows.Item(2).RowHeight = 17

# Update workbook window horizontal position (cosmetic, from saved window state)
$wb.Windows.Item(1).Left = 7060
